# Edit the "conclusion" slide's repository-link textbox:
#   - left-align both paragraphs explicitly
#   - point the link text at the new fork/URL (with trailing ".git")
#   - resize the (spAutoFit) textbox to the new, slightly wider/shorter box

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(9)                 # "conclusion" slide
$sh = $s.Shapes.Item(3)                 # "TextBox 2" - the repo-link box
$tf = $sh.TextFrame
$tr = $tf.TextRange

# --- Paragraph 1: "You can find the code in this repository:" ---
$para1 = $tr.Paragraphs(1, 1)
# Force PowerPoint to materialise an explicit <a:pPr algn="l"/> by actually
# changing the alignment value before setting it back to left.
$para1.ParagraphFormat.Alignment = 2   # ppAlignCenter (dummy, forces a diff)
$para1.ParagraphFormat.Alignment = 1   # ppAlignLeft

# --- Paragraph 2: the hyperlinked repository URL ---
$para2 = $tr.Paragraphs(2, 1)
$para2.ParagraphFormat.Alignment = 2   # ppAlignCenter (dummy, forces a diff)
$para2.ParagraphFormat.Alignment = 1   # ppAlignLeft

# Replace just the run's characters (keeps the existing run/formatting,
# including the hyperlink + accent colour, instead of fragmenting runs).
$para2 = $tr.Paragraphs(2, 1)
$runChars = $para2.Characters(1, $para2.Length)
$runChars.Text = "https://github.com/arijitgolui10/VaultofCodes_Python_Intern_FinalProject.git"

# --- Resize the textbox (it uses <a:spAutoFit/>, so the author's PowerPoint
# reflowed it when the text changed) ---
$sh.Width  = 698.4
$sh.Height = 50.80001
